# Update odds values in row 4 of the active worksheet to match the
# FlashScore data refresh described in the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G4"  = 2.57
    "H4"  = 3.1
    "L4"  = 3.3
    "T4"  = 2.45
    "V4"  = 1.98
    "W4"  = 8.75
    "AB4" = 28
    "AC4" = 9.25
    "AH4" = 8.5
    "AI4" = 13.5
    "AJ4" = 9.75
    "AK4" = 30
    "AM4" = 30
    "AN4" = 4.45
    "AQ4" = 60
    "AR4" = 90
    "AT4" = 2.42
    "AU4" = 6.8
    "AY4" = 23
    "AZ4" = 70
    "BA4" = 110
    "BB4" = 300
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
